$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 39,40,42,43,44,45,46,47,48 hold the same 9 data records as before,
# just relocated to different row numbers (row 41 is untouched). Column I
# (Antal) is stored as text in this sheet, so values are written with a
# leading apostrophe to force text and avoid being coerced to numbers.

# --- Row 39 ---
$ws.Range("A39").Value = 111880601
$ws.Range("B39").Value = 88966
$ws.Range("E39").Value = 5754
$ws.Range("F39").Value = "Gultoppig fingersvamp"
$ws.Range("G39").Value = "Ramaria testaceoflava"
$ws.Range("H39").Value = "(Bres.) Corner"
$ws.Range("I39").Value = "'4"
$ws.Range("Q39").Value = 509941.5744066621
$ws.Range("R39").Value = 6753224.672924293

# --- Row 40 ---
$ws.Range("A40").Value = 111880475
$ws.Range("B40").Value = 88966
$ws.Range("E40").Value = 5754
$ws.Range("F40").Value = "Gultoppig fingersvamp"
$ws.Range("G40").Value = "Ramaria testaceoflava"
$ws.Range("H40").Value = "(Bres.) Corner"
$ws.Range("I40").Value = "'2"
$ws.Range("Q40").Value = 509957.7514087428
$ws.Range("R40").Value = 6753362.853637428
$ws.Range("AJ40").Value = "gran"
$ws.Range("AK40").Value = "Picea abies"
$ws.Range("AO40").Value = "Picea abies"

# --- Row 42 ---
$ws.Range("A42").Value = 111880484
$ws.Range("I42").Value = "'11"
$ws.Range("Q42").Value = 509900.7891887496
$ws.Range("R42").Value = 6753525.142772059

# --- Row 43 ---
$ws.Range("A43").Value = 111880574
$ws.Range("I43").Value = "'2"
$ws.Range("Q43").Value = 509595.7160662179
$ws.Range("R43").Value = 6753391.52735021

# --- Row 44 ---
$ws.Range("A44").Value = 111880462
$ws.Range("B44").Value = 88966
$ws.Range("E44").Value = 5754
$ws.Range("F44").Value = "Gultoppig fingersvamp"
$ws.Range("G44").Value = "Ramaria testaceoflava"
$ws.Range("H44").Value = "(Bres.) Corner"
$ws.Range("I44").Value = "'1"
$ws.Range("Q44").Value = 509970.2466718731
$ws.Range("R44").Value = 6753250.046013334
$ws.Range("AL44").Value = "vid tallar"
$ws.Range("AO44").Value = "Pinus sylvestris # vid tallar"

# --- Row 45 ---
$ws.Range("A45").Value = 111880591
$ws.Range("B45").Value = 90658
$ws.Range("E45").Value = 4361
$ws.Range("F45").Value = "Orange taggsvamp"
$ws.Range("G45").Value = "Hydnellum aurantiacum"
$ws.Range("H45").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("I45").Value = "'8"
$ws.Range("Q45").Value = 509822.1902239832
$ws.Range("R45").Value = 6753234.069152902
$ws.Range("AJ45").Value = "tall"
$ws.Range("AK45").Value = "Pinus sylvestris"
$ws.Range("AO45").Value = "Pinus sylvestris"

# --- Row 46 ---
$ws.Range("A46").Value = 111880509
$ws.Range("B46").Value = 90652
$ws.Range("E46").Value = 3100
$ws.Range("F46").Value = "Talltaggsvamp"
$ws.Range("G46").Value = "Bankera fuligineoalba"
$ws.Range("H46").Value = "(Schmidt : Fr.) Pouzar"
$ws.Range("I46").Value = "'6"
$ws.Range("Q46").Value = 509834.2096935506
$ws.Range("R46").Value = 6753644.114383955

# --- Row 47 ---
$ws.Range("A47").Value = 111880562
$ws.Range("B47").Value = 90658
$ws.Range("E47").Value = 4361
$ws.Range("F47").Value = "Orange taggsvamp"
$ws.Range("G47").Value = "Hydnellum aurantiacum"
$ws.Range("H47").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("I47").Value = "'3"
$ws.Range("Q47").Value = 509657.7198006394
$ws.Range("R47").Value = 6753521.069647122

# --- Row 48 ---
$ws.Range("A48").Value = 111880580
$ws.Range("B48").Value = 90658
$ws.Range("E48").Value = 4361
$ws.Range("F48").Value = "Orange taggsvamp"
$ws.Range("G48").Value = "Hydnellum aurantiacum"
$ws.Range("H48").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("I48").Value = "'3"
$ws.Range("Q48").Value = 509755.441071702
$ws.Range("R48").Value = 6753236.317390828
$ws.Range("AO48").Value = "Pinus sylvestris"
$ws.Range("AL48").ClearContents()
